# Updated symbol list on Sun Dec 18 10:46:24 UTC 2022 with GitHub Actions
#
# Applies the refreshed "Price" (column D) quotes and the handful of
# "Volume(1h)" label (column E) / row-ordering corrections that came out
# of that run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    # Force text storage so numeric-looking strings (e.g. "247.07") are not
    # reinterpreted as numbers by Excel's auto-detection.
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    # Revert the cell back to the workbook's default style so no stray
    # per-cell formatting is introduced.
    $rng.Style = "Normal"
}

# --- Column D ("Price") quote refreshes -------------------------------
Set-TextValue "D2"  "247.07"
Set-TextValue "D3"  "22.40"
Set-TextValue "D4"  "5.478"
Set-TextValue "D5"  "0.05610"
Set-TextValue "D6"  "6.468"
Set-TextValue "D7"  "0.8040"
Set-TextValue "D9"  "0.1424"
Set-TextValue "D10" "0.07284"
Set-TextValue "D11" "0.03177"
Set-TextValue "D12" "0.02948"
Set-TextValue "D13" "0.09260"
Set-TextValue "D15" "3.196"
Set-TextValue "D16" "0.04699"
Set-TextValue "D17" "0.0005983"
Set-TextValue "D18" "0.006276"
Set-TextValue "D19" "0.001060"
Set-TextValue "D20" "0.003817"
Set-TextValue "D21" "0.0001503"
Set-TextValue "D22" "0.0003307"
Set-TextValue "D23" "3.980"
Set-TextValue "D24" "3.393"
Set-TextValue "D41" "0.1044"
Set-TextValue "D42" "0.002976"
Set-TextValue "D43" "0.003253"
Set-TextValue "D44" "0.01025"
Set-TextValue "D45" "0.00005632"
Set-TextValue "D47" "0.6815"
Set-TextValue "D48" "0.02599"
Set-TextValue "D49" "0.00002104"

# --- Column E ("Volume(1h)") label refreshes ---------------------------
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E27").Value = "26ProBitTokenPROBBestin24h"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

# --- Rows 41 / 43 swap: KickToken <-> BKEXToken -------------------------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
